# Update the metrics sheet: column A model names get reshuffled across rows
# 2-26, and every metric column (B:Q) for those rows is overwritten with a
# single shared set of recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-name order for rows 2..26 (row r -> $names[r-2])
$names = @(
    "model_6_5_0",
    "model_6_5_22",
    "model_6_5_21",
    "model_6_5_20",
    "model_6_5_19",
    "model_6_5_18",
    "model_6_5_17",
    "model_6_5_16",
    "model_6_5_15",
    "model_6_5_14",
    "model_6_5_13",
    "model_6_5_23",
    "model_6_5_12",
    "model_6_5_10",
    "model_6_5_9",
    "model_6_5_8",
    "model_6_5_7",
    "model_6_5_6",
    "model_6_5_5",
    "model_6_5_4",
    "model_6_5_3",
    "model_6_5_2",
    "model_6_5_1",
    "model_6_5_11",
    "model_6_5_24"
)

# Shared metric values (columns B..Q) now common to every data row (2..26)
$metricCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")
$metricVals = @(
    0.5692387035740197,
    0.2355707868591003,
    -2.163366163360928,
    -16.99369045119845,
    -0.7920686400455674,
    0.2557185134775685,
    0.4537982025430372,
    0.5733796306930626,
    0.152569633957454,
    0.3629746323252583,
    0.2823828768330753,
    0.5056861808251918,
    0.06015717143422494,
    0.5272142884576241,
    28.72735599237364,
    44.57274171566024
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $metricCols.Length; $j++) {
        $ws.Range($metricCols[$j] + $row).Value = $metricVals[$j]
    }
}
